$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.63"
$ws.Range("E2").Value = "'0.71%"
$ws.Range("D3").Value = "'29.82"
$ws.Range("E3").Value = "'10.01%"
$ws.Range("D4").Value = "'5.181"
$ws.Range("E4").Value = "'1.83%"
$ws.Range("D5").Value = "'0.05726"
$ws.Range("E5").Value = "'0.43%"
$ws.Range("D6").Value = "'6.598"
$ws.Range("E6").Value = "'1.42%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.087"
$ws.Range("E7").Value = "'2.59%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8566"
$ws.Range("E8").Value = "'4.43%"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8704"
$ws.Range("E9").Value = "'0.84%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1361"
$ws.Range("E10").Value = "'2.21%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07039"
$ws.Range("E11").Value = "'1.92%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02924"
$ws.Range("E12").Value = "'3.41%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09386"
$ws.Range("E13").Value = "'0.15%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001527"
$ws.Range("E14").Value = "'0.63%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04173"
$ws.Range("E15").Value = "'2.95%"
$ws.Range("D16").Value = "'0.006084"
$ws.Range("E16").Value = "'0.72%"
$ws.Range("E17").Value = "'5,069.54%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.269"
$ws.Range("E19").Value = "'-1.99%"
$ws.Range("B20").Value = "One"
$ws.Range("C20").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D20").Value = "'0.01028"
$ws.Range("E20").Value = "'1,613.05%"
$ws.Range("D21").Value = "'0.3156"
$ws.Range("E21").Value = "'-0.68%"
$ws.Range("D22").Value = "'0.03381"
$ws.Range("E22").Value = "'6.55%"
$ws.Range("D23").Value = "'0.1316"
$ws.Range("E23").Value = "'3.29%"
$ws.Range("D24").Value = "'3.464"
$ws.Range("E24").Value = "'-2.84%"
$ws.Range("D25").Value = "'0.1379"
$ws.Range("E25").Value = "'0.40%"
$ws.Range("D26").Value = "'0.005023"
$ws.Range("E26").Value = "'26.44%"
$ws.Range("E27").Value = "'0.34%"
$ws.Range("D28").Value = "'0.0001209"
$ws.Range("E28").Value = "'22.17%"
$ws.Range("D40").Value = "'0.03752"
$ws.Range("E40").Value = "'0.82%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005747"
$ws.Range("E41").Value = "'0.66%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E42").Value = "'1.41%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.001999"
$ws.Range("E43").Value = "'-15.58%"
$ws.Range("D44").Value = "'0.009570"
$ws.Range("E44").Value = "'2.08%"
$ws.Range("D45").Value = "'0.00005224"
$ws.Range("E45").Value = "'1.05%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("D47").Value = "'0.06466"
$ws.Range("E47").Value = "'-36.28%"
$ws.Range("D48").Value = "'0.002517"
$ws.Range("E48").Value = "'-0.79%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.03%"
